$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old "Terms Typically Offered" column (D),
# shifting its data to column G.
$ws.Columns("D:F").Insert()

# New header row values
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default new columns to "NA" for all data rows (2-20)
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 4).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
}

# Row 9: split Corequisite info out of Prerequisites (column C)
$ws.Range("C9").Value = "Junior standing."
$ws.Range("D9").Value = "LS 201."
$ws.Range("G9").Value = "F, SP "

# Row 12: split Recommended info out of Prerequisites (column C)
$ws.Range("C12").Value = "LS 211 and GE C3."
$ws.Range("F12").Value = "LS 310."
$ws.Range("G12").Value = "F, W, SP "

# Row 17: split Recommended info out of Prerequisites (column C)
$ws.Range("C17").Value = "PSC 103; BIO 211; LS 250; MATH 328; and junior standing."
$ws.Range("F17").Value = "MATH 329."
$ws.Range("G17").Value = "F, SP "
